# "connected robotrainer and program creator"
# Extends the existing 4-day workout program (WEEK 1) out to 9 days,
# and refreshes the exercise/rep values for every day so they match
# the values produced by the program-creator tool.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Full target data for every day block: header text + the 4 exercises
# (name, reps) in order (Cable Flies, DB Press, Incline DB Press,
# Incline Press Machine).
$days = @(
    @{ Header = "DAY 1"; Exercises = @(
            @("Cable Flies", "4x6"),
            @("DB Press", "3x12"),
            @("Incline DB Press", "3x10"),
            @("Incline Press Machine", "4x8")
        ) },
    @{ Header = "DAY 2"; Exercises = @(
            @("Cable Flies", "3x12"),
            @("DB Press", "4x10"),
            @("Incline DB Press", "3x6"),
            @("Incline Press Machine", "4x8")
        ) },
    @{ Header = "DAY 3"; Exercises = @(
            @("Cable Flies", "4x12"),
            @("DB Press", "4x10"),
            @("Incline DB Press", "4x10"),
            @("Incline Press Machine", "3x6")
        ) },
    @{ Header = "DAY 4"; Exercises = @(
            @("Cable Flies", "4x10"),
            @("DB Press", "3x12"),
            @("Incline DB Press", "4x10"),
            @("Incline Press Machine", "3x12")
        ) },
    @{ Header = "DAY 5"; Exercises = @(
            @("Cable Flies", "4x8"),
            @("DB Press", "4x12"),
            @("Incline DB Press", "3x12"),
            @("Incline Press Machine", "3x12")
        ) },
    @{ Header = "DAY 6"; Exercises = @(
            @("Cable Flies", "3x6"),
            @("DB Press", "4x6"),
            @("Incline DB Press", "3x10"),
            @("Incline Press Machine", "3x8")
        ) },
    @{ Header = "DAY 7"; Exercises = @(
            @("Cable Flies", "4x12"),
            @("DB Press", "4x6"),
            @("Incline DB Press", "3x8"),
            @("Incline Press Machine", "3x6")
        ) },
    @{ Header = "DAY 8"; Exercises = @(
            @("Cable Flies", "3x8"),
            @("DB Press", "3x10"),
            @("Incline DB Press", "3x6"),
            @("Incline Press Machine", "3x8")
        ) },
    @{ Header = "DAY 9"; Exercises = @(
            @("Cable Flies", "3x10"),
            @("DB Press", "3x12"),
            @("Incline DB Press", "4x10"),
            @("Incline Press Machine", "3x12")
        ) }
)

# Each day block occupies 8 rows: header row, blank row, 4 exercise
# rows, then 2 blank rows before the next header. Day 1 starts at
# row 5, so day N's header is at row (5 + 8*(N-1)).
$rowsPerBlock = 8
$firstHeaderRow = 5

# Only days 5-9 are brand new rows; days 1-4 (rows 5, 13, 21, 29) already
# exist in the sheet and must only get their values refreshed below.
$existingDayCount = 4
$newHeaderRows = @()
for ($i = $existingDayCount; $i -lt $days.Count; $i++) {
    $newHeaderRows += ($firstHeaderRow + ($rowsPerBlock * $i))
}

# Days 5-9 (header rows 37, 45, 53, 61, 69) don't exist yet in the
# original sheet. Recreate the merged cells for all of them first, then
# clone the formatting from the Day 1 block (header + 4 exercise rows)
# on top of the merges. Doing all the merges before any of the
# formatting paste keeps the style table from forking (merging an
# already-formatted/pasted range duplicates styles in this engine).
foreach ($headerRow in $newHeaderRows) {
    $firstExerciseRow = $headerRow + 2
    $ws.Range("A" + $headerRow + ":B" + $headerRow).Merge() | Out-Null
    for ($e = 0; $e -lt 4; $e++) {
        $r = $firstExerciseRow + $e
        $ws.Range("A" + $r + ":C" + $r).Merge() | Out-Null
    }
}

foreach ($headerRow in $newHeaderRows) {
    $ws.Range("A5:B5").Copy() | Out-Null
    $ws.Range("A$headerRow").PasteSpecial($xlPasteFormats) | Out-Null
}
foreach ($headerRow in $newHeaderRows) {
    $firstExerciseRow = $headerRow + 2
    $ws.Range("A7:D10").Copy() | Out-Null
    $ws.Range("A$firstExerciseRow").PasteSpecial($xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $days.Count; $i++) {
    $day = $days[$i]
    $headerRow = $firstHeaderRow + ($rowsPerBlock * $i)
    $firstExerciseRow = $headerRow + 2

    $ws.Range("A$headerRow").Value = $day.Header

    for ($e = 0; $e -lt $day.Exercises.Count; $e++) {
        $r = $firstExerciseRow + $e
        $pair = $day.Exercises[$e]
        $ws.Range("A$r").Value = $pair[0]
        $ws.Range("D$r").Value = $pair[1]
    }
}
